$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 44179
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 3000

# Row 5 updates
$ws.Range("D5").Value = 45243
$ws.Range("M5").Value = 50

# Row 6 updates
$ws.Range("D6").Value = 45250
$ws.Range("M6").Value = 30
$ws.Range("R6").Value = "Provincia de San Felipe de Aconcagua"

# Row 7 updates
$ws.Range("D7").Value = 45244
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 35000
$ws.Range("O7").Value = 35000
$ws.Range("P7").Value = 35000
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 7000
